$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$recordedBy = @{
    2 = "eman.tantawi@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg"
    3 = "eman.tantawi@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
    4 = "eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
    5 = "nourhan.mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, hananragab@med.asu.edu.eg"
    6 = "nourhan.mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg"
    7 = "eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
    8 = "eman.tantawi@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, manar.montaser@med.asu.edu.eg, backup@backdoor.com, asmaa.reda@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"
    9 = "hend_mahmoud@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, gehanadel@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"
    13 = "mariam.noureldin@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg, wessam.atef@med.asu.edu.eg"
    15 = "aml.awwad@med.asu.edu.eg, wessam.atef@med.asu.edu.eg"
    16 = "aml.awwad@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg"
    24 = "Salma.hassan@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, marina_atef@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg"
    25 = "Aya_hamed@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, marina_atef@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"
    27 = "Salma.hassan@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Eman_mohamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg"
    28 = "Salma.hassan@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"
    30 = "eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
    31 = "eman.tantawi@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
    32 = "eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
    33 = "nourhan.mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, hananragab@med.asu.edu.eg"
    34 = "nourhan.mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg"
    35 = "eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
    36 = "eman.tantawi@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, manar.montaser@med.asu.edu.eg, backup@backdoor.com, asmaa.reda@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"
    37 = "hend_mahmoud@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, gehanadel@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"
    41 = "mariam.noureldin@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg, wessam.atef@med.asu.edu.eg"
    43 = "aml.awwad@med.asu.edu.eg, wessam.atef@med.asu.edu.eg"
    44 = "aml.awwad@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg"
    52 = "Salma.hassan@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, marina_atef@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg"
    53 = "Aya_hamed@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, marina_atef@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"
    55 = "Salma.hassan@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Eman_mohamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg"
    56 = "Salma.hassan@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"
}

foreach ($row in $recordedBy.Keys) {
    $ws.Cells.Item($row, 7).Value = $recordedBy[$row]
}

